$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Horario")
$ws2 = $wb.Worksheets.Item("Detalle_horario")

# Insert a new column before column A on both sheets, shifting existing data right.
$ws1.Columns.Item(1).Insert()
$ws2.Columns.Item(1).Insert()

# Set the new header cell (column A, row 1) to "N" on both sheets.
$ws1.Range("A1").Value = "N"
$ws2.Range("A1").Value = "N"
